$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.006.34'
$ws.Range('E2').Value = '  +0.67%  '

$ws.Range('D3').Value = '1.641.79'
$ws.Range('E3').Value = '  +0.75%  '

$ws.Range('E4').Value = '  +0.33%  '

$ws.Range('D5').Value = '215.90'
$ws.Range('E5').Value = '  +0.93%  '

$ws.Range('E6').Value = '  +0.27%  '

$ws.Range('E7').Value = '  +0.31%  '

$ws.Range('E8').Value = '  +0.48%  '

$ws.Range('E9').Value = '  +1.09%  '

$ws.Range('D10').Value = '19.57'
$ws.Range('E10').Value = '  +0.36%  '

$ws.Range('D11').Value = '0.0796'
$ws.Range('E11').Value = '  +0.56%  '

$ws.Range('D12').Value = '1.869.51'
$ws.Range('E12').Value = '  +0.72%  '

$ws.Range('E13').Value = '  +0.76%  '

$ws.Range('D14').Value = '1.647.46'
$ws.Range('E14').Value = '  +1.22%  '

$ws.Range('E15').Value = '  +0.24%  '

$ws.Range('D16').Value = '0.0₃0763'
$ws.Range('E16').Value = '  +1.05%  '

$ws.Range('E17').Value = '  +1.35%  '

$ws.Range('D18').Value = '26.105.20'

$ws.Range('E19').Value = '  +0.39%  '

$ws.Range('D20').Value = '194.42'
$ws.Range('E20').Value = '  +0.73%  '

$ws.Range('E21').Value = '  -0.47%  '

$ws.Range('E22').Value = '  +0.29%  '

$ws.Range('D23').Value = '6.21'
$ws.Range('E23').Value = '  -0.74%  '

$ws.Range('D24').Value = '1.81'
$ws.Range('E24').Value = '  -0.63%  '

$ws.Range('E25').Value = '  +4.75%  '

$ws.Range('E26').Value = '  +0.41%  '

$ws.Range('D27').Value = '143.01'
$ws.Range('E27').Value = '  -0.15%  '

$ws.Range('E28').Value = '  +0.79%  '

$ws.Range('E29').Value = '  +0.86%  '

$ws.Range('E30').Value = '  +0.72%  '

$ws.Range('D31').Value = '0.0496'
$ws.Range('E31').Value = '  -0.53%  '

$ws.Range('E32').Value = '  -0.14%  '

$ws.Range('E33').Value = '  +1.26%  '

$ws.Range('E34').Value = '  -0.88%  '

$ws.Range('E35').Value = '  +1.27%  '

$ws.Range('D36').Value = '0.905'
$ws.Range('E36').Value = '  +0.55%  '

$ws.Range('D37').Value = '1.130.56'
$ws.Range('E37').Value = '  -0.57%  '

$ws.Range('D38').Value = '0.540'
$ws.Range('E38').Value = '  -1.00%  '

$ws.Range('E39').Value = '  -0.48%  '

$ws.Range('E40').Value = '  +0.27%  '

$ws.Range('E41').Value = '  +0.93%  '

$ws.Range('D42').Value = '99.08'
$ws.Range('E42').Value = '  -0.03%  '

$ws.Range('E43').Value = '  +0.11%  '

$ws.Range('D44').Value = '1.778.21'
$ws.Range('E44').Value = '  +0.71%  '

$ws.Range('E45').Value = '  +3.93%  '

$ws.Range('D46').Value = '56.47'
$ws.Range('E46').Value = '  +0.58%  '

$ws.Range('E47').Value = '  -1.25%  '

$ws.Range('E48').Value = '  +3.26%  '

$ws.Range('D49').Value = '7.77'
$ws.Range('E49').Value = '  +1.59%  '

$ws.Range('E50').Value = '  -0.17%  '

$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.17%  '
